$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "葡萄糖"
$ws.Range("A3").Value = "尿胆原"
$ws.Range("A4").Value = "胆红素"
$ws.Range("A5").Value = "酮体"
$ws.Range("A6").Value = "比重"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1.02"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = "1.003-1.030"
$ws.Range("A7").Value = "酸碱度"
$ws.Range("C7").Value = "4.5-8"
$ws.Range("A8").Value = "蛋白质"
$ws.Range("A9").Value = "亚硝酸盐"
$ws.Range("A10").Value = "白细胞"
$ws.Range("A11").Value = "潜血"
$ws.Range("A12").Value = "RDW-CV"
$ws.Range("C12").Value = "0-5"
$ws.Range("A13").Value = "RDW-CV"
$ws.Range("C13").Value = "0-7"
$ws.Range("A14").Value = "上皮细胞"
$ws.Range("C14").Value = "0-5"
$ws.Range("A15").Value = "透明管型"
$ws.Range("C15").Value = "0-1"
$ws.Range("A16").Value = "细胞管型"
$ws.Range("A17").Value = "颗粒管型"
$ws.Range("A18").Value = "其他管型"
$ws.Range("C18").Value = "0-0"
$ws.Range("A19").Value = "尿酸盐结晶"
$ws.Range("A20").Value = "其他结晶"
$ws.Range("A21").Value = "草酸钙结晶"
$ws.Range("C21").Value = "0-0"
$ws.Range("A22").Value = "RDW-CV"
$ws.Range("A23").Value = "RDW-CV"
$ws.Range("C23").Value = "0-0"
